$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "time_taken" column, formatted like the other headers
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Timestamps for each data row (rows 2-54), matching values from the source commit
$timeTaken = @(
    "2021-10-05 13:39:07.194445",
    "2021-10-05 13:39:07.194459",
    "2021-10-05 13:39:07.194464",
    "2021-10-05 13:39:07.194467",
    "2021-10-05 13:39:07.194471",
    "2021-10-05 13:39:07.194474",
    "2021-10-05 13:39:07.194477",
    "2021-10-05 13:39:07.194481",
    "2021-10-05 13:39:07.194484",
    "2021-10-05 13:39:07.194487",
    "2021-10-05 13:39:07.194491",
    "2021-10-05 13:39:07.194494",
    "2021-10-05 13:39:07.194497",
    "2021-10-05 13:39:07.194500",
    "2021-10-05 13:39:07.194503",
    "2021-10-05 13:39:07.194506",
    "2021-10-05 13:39:07.194509",
    "2021-10-05 13:39:07.194512",
    "2021-10-05 13:39:07.194516",
    "2021-10-05 13:39:07.194519",
    "2021-10-05 13:39:07.194522",
    "2021-10-05 13:39:07.194525",
    "2021-10-05 13:39:07.194528",
    "2021-10-05 13:39:07.194531",
    "2021-10-05 13:39:07.194534",
    "2021-10-05 13:39:07.194538",
    "2021-10-05 13:39:07.194541",
    "2021-10-05 13:39:07.194544",
    "2021-10-05 13:39:07.194547",
    "2021-10-05 13:39:07.194550",
    "2021-10-05 13:39:07.194553",
    "2021-10-05 13:39:07.194556",
    "2021-10-05 13:39:07.194560",
    "2021-10-05 13:39:07.194564",
    "2021-10-05 13:39:07.194567",
    "2021-10-05 13:39:07.194570",
    "2021-10-05 13:39:07.194573",
    "2021-10-05 13:39:07.194576",
    "2021-10-05 13:39:07.194579",
    "2021-10-05 13:39:07.194583",
    "2021-10-05 13:39:07.194587",
    "2021-10-05 13:39:07.194590",
    "2021-10-05 13:39:07.194593",
    "2021-10-05 13:39:07.194596",
    "2021-10-05 13:39:07.194600",
    "2021-10-05 13:39:07.194603",
    "2021-10-05 13:39:07.194606",
    "2021-10-05 13:39:07.194609",
    "2021-10-05 13:39:07.194612",
    "2021-10-05 13:39:07.194615",
    "2021-10-05 13:39:07.194618",
    "2021-10-05 13:39:07.194621",
    "2021-10-05 13:39:07.194625"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "Added time_taken column with $($timeTaken.Length) values"
